# Update the "last_edited_time" column (D) for the Notion data refresh.
# Rows 4 and 6 move to the next edit timestamp (03:28); rows 5, 7, 8, 12 and
# 13 move to an even later timestamp (03:29) that did not previously exist
# in the workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D4").Value = "2024-08-03T03:28:00.000Z"
$ws.Range("D6").Value = "2024-08-03T03:28:00.000Z"

$ws.Range("D5").Value = "2024-08-03T03:29:00.000Z"
$ws.Range("D7").Value = "2024-08-03T03:29:00.000Z"
$ws.Range("D8").Value = "2024-08-03T03:29:00.000Z"
$ws.Range("D12").Value = "2024-08-03T03:29:00.000Z"
$ws.Range("D13").Value = "2024-08-03T03:29:00.000Z"
